# Auto-generated edit script applying scheduled-runner price/profit updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 7777.5
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 7777.5
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 7777.5
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -8589.5

$ws.Range("H91").Value = 7777.5
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 7777.5
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 7777.5
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -10585.5

$ws.Range("H98").Value = 3833.25
$ws.Range("I98").Value = 3833.25
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 3833.25
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -2335.25
$ws.Range("N98").ClearContents()

$ws.Range("H116").Value = 12215.333
$ws.Range("J116").Value = 6171.6924
$ws.Range("L116").Value = 6171.6924
$ws.Range("N116").Value = -13055.6924

$ws.Range("H122").Value = 3833.25
$ws.Range("I122").Value = 3833.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11499.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -9049.75
$ws.Range("N122").ClearContents()

$ws.Range("H137").Value = 1790.05
$ws.Range("I137").Value = 985.7857
$ws.Range("J137").Value = 3666.6667
$ws.Range("K137").Value = 2957.3571
$ws.Range("L137").Value = 11000.0001
$ws.Range("M137").Value = -407.3571000000002
$ws.Range("N137").Value = -16100.0001

$ws.Range("H138").Value = 4448.8335
$ws.Range("J138").Value = 3099.4707
$ws.Range("L138").Value = 9298.4121
$ws.Range("N138").Value = -19578.4121

$ws.Range("H141").Value = 1219109
$ws.Range("I141").Value = 1556453.2
$ws.Range("J141").Value = 4670
$ws.Range("K141").Value = 4669359.6
$ws.Range("L141").Value = 14010
$ws.Range("M141").Value = -4664179.6
$ws.Range("N141").Value = -24370

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2980.309
$ws.Range("I32").Value = 2752.7334
$ws.Range("J32").Value = 4004.4
$ws.Range("K32").Value = 2752.7334
$ws.Range("L32").Value = 4004.4
$ws.Range("M32").Value = -2465.7334
$ws.Range("N32").Value = -4578.4

$ws.Range("H45").Value = 1434.1724
$ws.Range("I45").Value = 1250.0555
$ws.Range("K45").Value = 1250.0555
$ws.Range("M45").Value = -873.0554999999999

$ws.Range("H61").Value = 2126.0698
$ws.Range("I61").Value = 1353.4736
$ws.Range("J61").Value = 7997.8
$ws.Range("K61").Value = 1353.4736
$ws.Range("L61").Value = 7997.8
$ws.Range("M61").Value = -1141.4736
$ws.Range("N61").Value = -8421.799999999999

$ws.Range("H122").Value = 1843.8125
$ws.Range("I122").Value = 1966.4615
$ws.Range("J122").Value = 1312.3334
$ws.Range("K122").Value = 5899.3845
$ws.Range("L122").Value = 3937.0002
$ws.Range("M122").Value = -3449.3845
$ws.Range("N122").Value = -8837.0002

$ws.Range("H132").Value = 1433.6364
$ws.Range("J132").Value = 3876.25
$ws.Range("L132").Value = 11628.75
$ws.Range("N132").Value = -16688.75

$ws.Range("H136").Value = 2126.0698
$ws.Range("I136").Value = 1353.4736
$ws.Range("J136").Value = 7997.8
$ws.Range("K136").Value = 4060.4208
$ws.Range("L136").Value = 23993.4
$ws.Range("M136").Value = -1510.4208
$ws.Range("N136").Value = -29093.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 14757
$ws.Range("J80").Value = 20548.4
$ws.Range("L80").Value = 20548.4
$ws.Range("N80").Value = -22544.4

$ws.Range("H83").Value = 14757
$ws.Range("J83").Value = 20548.4
$ws.Range("L83").Value = 102742
$ws.Range("N83").Value = -112726

$ws.Range("H134").Value = 4735.727
$ws.Range("I134").Value = 5284.1177
$ws.Range("K134").Value = 15852.3531
$ws.Range("M134").Value = -13317.3531

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1781.9375
$ws.Range("I31").Value = 1779.2
$ws.Range("K31").Value = 1779.2
$ws.Range("M31").Value = -1484.2

$ws.Range("H32").Value = 3010
$ws.Range("I32").Value = 3010
$ws.Range("K32").Value = 3010
$ws.Range("M32").Value = -2694

$ws.Range("H34").Value = 1781.9375
$ws.Range("I34").Value = 1779.2
$ws.Range("K34").Value = 1779.2
$ws.Range("M34").Value = -1577.2

$ws.Range("H58").Value = 1243437
$ws.Range("I58").Value = 1673110.9
$ws.Range("J58").Value = 2156.889
$ws.Range("K58").Value = 1673110.9
$ws.Range("L58").Value = 2156.889
$ws.Range("M58").Value = -1672907.9
$ws.Range("N58").Value = -2562.889

$ws.Range("H122").Value = 3495.75
$ws.Range("I122").Value = 1680.2222
$ws.Range("J122").Value = 8942.333000000001
$ws.Range("K122").Value = 5040.6666
$ws.Range("L122").Value = 26826.999
$ws.Range("M122").Value = -2590.6666
$ws.Range("N122").Value = -31726.999

$ws.Range("H132").Value = 1532.2
$ws.Range("I132").Value = 1018.75
$ws.Range("K132").Value = 3056.25
$ws.Range("M132").Value = -526.25

$ws.Range("H134").Value = 1141.5555
$ws.Range("I134").Value = 1059.8889
$ws.Range("J134").Value = 1468.2222
$ws.Range("K134").Value = 3179.6667
$ws.Range("L134").Value = 4404.6666
$ws.Range("M134").Value = -644.6666999999998
$ws.Range("N134").Value = -9474.6666

$ws.Range("H136").Value = 1243437
$ws.Range("I136").Value = 1673110.9
$ws.Range("J136").Value = 2156.889
$ws.Range("K136").Value = 5019332.699999999
$ws.Range("L136").Value = 6470.667
$ws.Range("M136").Value = -5016782.699999999
$ws.Range("N136").Value = -11570.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 2000
$ws.Range("I25").Value = 8000
$ws.Range("J25").Value = 1400
$ws.Range("K25").Value = 24000
$ws.Range("L25").Value = 4200
$ws.Range("M25").Value = -23831
$ws.Range("N25").Value = -4538

$ws.Range("H30").Value = 2000
$ws.Range("I30").Value = 8000
$ws.Range("J30").Value = 1400
$ws.Range("K30").Value = 24000
$ws.Range("L30").Value = 4200
$ws.Range("M30").Value = -23898
$ws.Range("N30").Value = -4404

$ws.Range("H107").Value = 891.05554
$ws.Range("I107").Value = 604.5714
$ws.Range("J107").Value = 1073.3636
$ws.Range("K107").Value = 1813.7142
$ws.Range("L107").Value = 3220.0908
$ws.Range("M107").Value = 106.2857999999999
$ws.Range("N107").Value = -7060.0908

$ws.Range("H122").Value = 948.3889
$ws.Range("J122").Value = 1156.9166
$ws.Range("L122").Value = 10412.2494
$ws.Range("N122").Value = -15312.2494

$ws.Range("H131").Value = 9187.393
$ws.Range("J131").Value = 9833.73
$ws.Range("L131").Value = 29501.19
$ws.Range("N131").Value = -39581.19

$ws.Range("H136").Value = 4000
$ws.Range("I136").Value = 4000
$ws.Range("K136").Value = 12000
$ws.Range("M136").Value = -6900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

$ws.Range("H35").Value = 30008
$ws.Range("I35").Value = 2999
$ws.Range("K35").Value = 2999
$ws.Range("M35").Value = -2701

$ws.Range("H36").Value = 1001
$ws.Range("I36").Value = 1001
$ws.Range("K36").Value = 1001
$ws.Range("M36").Value = -516

$ws.Range("H43").Value = 1900
$ws.Range("I43").Value = 1900
$ws.Range("K43").Value = 1900
$ws.Range("M43").Value = -1749

$ws.Range("H122").Value = 1489.3636
$ws.Range("I122").Value = 1067.7142
$ws.Range("J122").Value = 2227.25
$ws.Range("K122").Value = 3203.1426
$ws.Range("L122").Value = 6681.75
$ws.Range("M122").Value = -753.1425999999997
$ws.Range("N122").Value = -11581.75

$ws.Range("H126").Value = 2021954
$ws.Range("I126").Value = 4631656
$ws.Range("J126").Value = 64677.375
$ws.Range("K126").Value = 13894968
$ws.Range("L126").Value = 194032.125
$ws.Range("M126").Value = -13892498
$ws.Range("N126").Value = -198972.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1250
$ws.Range("J9").Value = 1250
$ws.Range("L9").Value = 1250
$ws.Range("N9").Value = -1698

$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

$ws.Range("H22").Value = 3344.4443
$ws.Range("J22").Value = 2585.7144
$ws.Range("L22").Value = 2585.7144
$ws.Range("N22").Value = -3175.7144

$ws.Range("H27").Value = 3344.4443
$ws.Range("J27").Value = 2585.7144
$ws.Range("L27").Value = 2585.7144
$ws.Range("N27").Value = -2799.7144

$ws.Range("H46").Value = 2142.1428
$ws.Range("I46").Value = 1748.75
$ws.Range("K46").Value = 1748.75
$ws.Range("M46").Value = -1560.75

$ws.Range("H61").Value = 1750.0834
$ws.Range("I61").Value = 1800.4667
$ws.Range("J61").Value = 1666.1111
$ws.Range("K61").Value = 1800.4667
$ws.Range("L61").Value = 1666.1111
$ws.Range("M61").Value = -1598.4667
$ws.Range("N61").Value = -2070.1111

$ws.Range("H113").Value = 1750.0834
$ws.Range("I113").Value = 1800.4667
$ws.Range("J113").Value = 1666.1111
$ws.Range("K113").Value = 1800.4667
$ws.Range("L113").Value = 1666.1111
$ws.Range("M113").Value = 369.5333000000001
$ws.Range("N113").Value = -6006.1111

$ws.Range("H132").Value = 1544.2963
$ws.Range("I132").Value = 1057.3572
$ws.Range("J132").Value = 2068.6924
$ws.Range("K132").Value = 3172.0716
$ws.Range("L132").Value = 6206.0772
$ws.Range("M132").Value = -642.0715999999998
$ws.Range("N132").Value = -11266.0772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 794.375
$ws.Range("I107").Value = 622.1429000000001
$ws.Range("K107").Value = 1866.4287
$ws.Range("M107").Value = 53.57129999999984

$ws.Range("H132").Value = 1761.1875
$ws.Range("J132").Value = 3992.3333
$ws.Range("L132").Value = 11976.9999
$ws.Range("N132").Value = -17036.9999

$ws.Range("H136").Value = 16837090
$ws.Range("I136").Value = 23150018
$ws.Range("K136").Value = 69450054
$ws.Range("M136").Value = -69447504

